$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in A2 and A3
$ws.Range("A2").Value = "39762 (non in estrazione)"
$ws.Range("A3").Value = "39666 (non in estrazione)"

# Add new row 4 with same pattern as rows 2 and 3
$ws.Range("A4").Value = "39742 (non in estrazione)"
$ws.Range("B4").Value = "CAMPO VUOTO"
$ws.Range("C4").Value = "CAMPO VUOTO"
$ws.Range("D4").Value = 0

# Copy the style (yellow fill) from B3/C3 to B4/C4
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
